$wb = $excel.ActiveWorkbook

# ALC row 28
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(28,8).Value = 45455204
$ws.Cells.Item(28,9).Value = 83333740
$ws.Cells.Item(28,10).Value = 956.1
$ws.Cells.Item(28,11).Value = 83333740
$ws.Cells.Item(28,12).Value = 956.1
$ws.Cells.Item(28,13).Value = -83333255
$ws.Cells.Item(28,14).Value = -1926.1

# ALC row 106
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(106,8).Value = 20837142
$ws.Cells.Item(106,9).Value = 25644046
$ws.Cells.Item(106,11).Value = 25644046
$ws.Cells.Item(106,13).Value = -25643415

# ALC row 138
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(138,8).Value = 1430.2188
$ws.Cells.Item(138,10).Value = 3566.6667
$ws.Cells.Item(138,12).Value = 10700.0001
$ws.Cells.Item(138,14).Value = -20980.0001

# ALC row 140
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(140,8).Value = 75697.16
$ws.Cells.Item(140,10).Value = 86821.71000000001
$ws.Cells.Item(140,12).Value = 86821.71000000001
$ws.Cells.Item(140,14).Value = -97181.71000000001

# ARM row 5
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(5,8).Value = 384.16666
$ws.Cells.Item(5,9).Value = 176.75
$ws.Cells.Item(5,10).Value = 799
$ws.Cells.Item(5,11).Value = 176.75
$ws.Cells.Item(5,12).Value = 799
$ws.Cells.Item(5,13).Value = -64.75
$ws.Cells.Item(5,14).Value = -1023

# ARM row 7
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(7,8).Value = 59999
$ws.Cells.Item(7,10).Value = 59999
$ws.Cells.Item(7,12).Value = 59999
$ws.Cells.Item(7,14).Value = -60227

# ARM row 38
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(38,8).Value = 9950
$ws.Cells.Item(38,9).Value = 0
$ws.Cells.Item(38,11).Value = 0
$ws.Cells.Item(38,13).ClearContents()

# ARM row 39
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(39,8).Value = 9994.75
$ws.Cells.Item(39,9).Value = 9994.75
$ws.Cells.Item(39,10).Value = 0
$ws.Cells.Item(39,11).Value = 9994.75
$ws.Cells.Item(39,12).Value = 0
$ws.Cells.Item(39,13).Value = -9474.75
$ws.Cells.Item(39,14).ClearContents()

# ARM row 110
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(110,8).Value = 2387.4092
$ws.Cells.Item(110,9).Value = 2358.238
$ws.Cells.Item(110,10).Value = 3000
$ws.Cells.Item(110,11).Value = 2358.238
$ws.Cells.Item(110,12).Value = 3000
$ws.Cells.Item(110,13).Value = -313.2379999999998
$ws.Cells.Item(110,14).Value = -7090

# BSM row 4
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(4,8).Value = 384.16666
$ws.Cells.Item(4,9).Value = 176.75
$ws.Cells.Item(4,10).Value = 799
$ws.Cells.Item(4,11).Value = 176.75
$ws.Cells.Item(4,12).Value = 799
$ws.Cells.Item(4,13).Value = -61.75
$ws.Cells.Item(4,14).Value = -1029

# BSM row 7
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(7,8).Value = 10501.333
$ws.Cells.Item(7,10).Value = 14502
$ws.Cells.Item(7,12).Value = 14502
$ws.Cells.Item(7,14).Value = -14728

# BSM row 38
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(38,8).Value = 0
$ws.Cells.Item(38,10).Value = 0
$ws.Cells.Item(38,12).Value = 0
$ws.Cells.Item(38,14).ClearContents()

# BSM row 105
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(105,8).Value = 5118.364
$ws.Cells.Item(105,9).Value = 3352.1333
$ws.Cells.Item(105,10).Value = 8903.143
$ws.Cells.Item(105,11).Value = 3352.1333
$ws.Cells.Item(105,12).Value = 8903.143
$ws.Cells.Item(105,13).Value = -1605.1333
$ws.Cells.Item(105,14).Value = -12397.143

# BSM row 107
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(107,8).Value = 4427.8105
$ws.Cells.Item(107,9).Value = 4146.0347
$ws.Cells.Item(107,10).Value = 4709.5864
$ws.Cells.Item(107,11).Value = 4146.0347
$ws.Cells.Item(107,12).Value = 4709.5864
$ws.Cells.Item(107,13).Value = -2226.0347
$ws.Cells.Item(107,14).Value = -8549.5864

# BSM row 116
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(116,8).Value = 110000
$ws.Cells.Item(116,10).Value = 110000
$ws.Cells.Item(116,12).Value = 110000
$ws.Cells.Item(116,14).Value = -119178

# CRP row 35
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(35,8).Value = 461.18182
$ws.Cells.Item(35,9).Value = 461.18182
$ws.Cells.Item(35,10).Value = 0
$ws.Cells.Item(35,11).Value = 461.18182
$ws.Cells.Item(35,12).Value = 0
$ws.Cells.Item(35,13).Value = -167.18182
$ws.Cells.Item(35,14).ClearContents()

# CRP row 64
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(64,8).Value = 42666
$ws.Cells.Item(64,10).Value = 58999.5
$ws.Cells.Item(64,12).Value = 58999.5
$ws.Cells.Item(64,14).Value = -59495.5

# CRP row 67
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(67,8).Value = 42666
$ws.Cells.Item(67,10).Value = 58999.5
$ws.Cells.Item(67,12).Value = 58999.5
$ws.Cells.Item(67,14).Value = -60715.5

# CRP row 122
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(122,8).Value = 2929.9824
$ws.Cells.Item(122,9).Value = 1562.2
$ws.Cells.Item(122,11).Value = 4686.6
$ws.Cells.Item(122,13).Value = -2236.6

# CUL row 5
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(5,8).Value = 699.2
$ws.Cells.Item(5,9).Value = 749
$ws.Cells.Item(5,11).Value = 2247
$ws.Cells.Item(5,13).Value = -2135

# CUL row 118
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(118,8).Value = 987
$ws.Cells.Item(118,9).Value = 987
$ws.Cells.Item(118,11).Value = 2961
$ws.Cells.Item(118,13).Value = -1718

# CUL row 122
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(122,8).Value = 1300.5454
$ws.Cells.Item(122,10).Value = 2499.75
$ws.Cells.Item(122,12).Value = 22497.75
$ws.Cells.Item(122,14).Value = -27397.75

# CUL row 131
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(131,8).Value = 1737493.9
$ws.Cells.Item(131,10).Value = 9048.615
$ws.Cells.Item(131,12).Value = 27145.845
$ws.Cells.Item(131,14).Value = -37225.845

# CUL row 132
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(132,8).Value = 2818.8823
$ws.Cells.Item(132,9).Value = 2478.2307
$ws.Cells.Item(132,11).Value = 22304.0763
$ws.Cells.Item(132,13).Value = -19774.0763

# CUL row 135
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(135,8).Value = 699.2
$ws.Cells.Item(135,9).Value = 749
$ws.Cells.Item(135,11).Value = 6741
$ws.Cells.Item(135,13).Value = -4206

# BSM row 45
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(45,8).Value = 169999
$ws.Cells.Item(45,10).Value = 169999
$ws.Cells.Item(45,12).Value = 169999
$ws.Cells.Item(45,14).Value = -171117

# GSM row 70
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(70,8).Value = 6971
$ws.Cells.Item(70,9).Value = 7798.8
$ws.Cells.Item(70,11).Value = 7798.8
$ws.Cells.Item(70,13).Value = -7528.8

# GSM row 73
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(73,8).Value = 6971
$ws.Cells.Item(73,9).Value = 7798.8
$ws.Cells.Item(73,11).Value = 7798.8
$ws.Cells.Item(73,13).Value = -6862.8

# GSM row 107
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(107,8).Value = 1021.4286
$ws.Cells.Item(107,9).Value = 838.5714
$ws.Cells.Item(107,11).Value = 838.5714
$ws.Cells.Item(107,13).Value = 1081.4286

# GSM row 122
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(122,8).Value = 4837.143
$ws.Cells.Item(122,9).Value = 5211.778
$ws.Cells.Item(122,10).Value = 4162.8
$ws.Cells.Item(122,11).Value = 15635.334
$ws.Cells.Item(122,12).Value = 12488.4
$ws.Cells.Item(122,13).Value = -13185.334
$ws.Cells.Item(122,14).Value = -17388.4

# LTW row 34
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(34,8).Value = 13000
$ws.Cells.Item(34,9).Value = 13000
$ws.Cells.Item(34,11).Value = 13000
$ws.Cells.Item(34,13).Value = -12828

# LTW row 48
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(48,8).Value = 29330
$ws.Cells.Item(48,9).Value = 29330
$ws.Cells.Item(48,11).Value = 29330
$ws.Cells.Item(48,13).Value = -28669

# LTW row 114
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(114,8).Value = 48900
$ws.Cells.Item(114,10).Value = 48900
$ws.Cells.Item(114,12).Value = 48900
$ws.Cells.Item(114,14).Value = -57578

# LTW row 122
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(122,8).Value = 17974.75
$ws.Cells.Item(122,9).Value = 17974.75
$ws.Cells.Item(122,11).Value = 53924.25
$ws.Cells.Item(122,13).Value = -51474.25

# LTW row 132
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(132,8).Value = 106106.09
$ws.Cells.Item(132,9).Value = 106106.09
$ws.Cells.Item(132,10).Value = 0
$ws.Cells.Item(132,11).Value = 318318.27
$ws.Cells.Item(132,12).Value = 0
$ws.Cells.Item(132,13).Value = -315788.27
$ws.Cells.Item(132,14).ClearContents()

# LTW row 136
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(136,8).Value = 10028.667
$ws.Cells.Item(136,9).Value = 8536.5
$ws.Cells.Item(136,11).Value = 25609.5
$ws.Cells.Item(136,13).Value = -23059.5

# WVR row 111
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(111,8).Value = 34500
$ws.Cells.Item(111,10).Value = 34500
$ws.Cells.Item(111,12).Value = 34500
$ws.Cells.Item(111,14).Value = -42680

# WVR row 132
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(132,8).Value = 3641.0857
$ws.Cells.Item(132,9).Value = 3719.4243
$ws.Cells.Item(132,10).Value = 2348.5
$ws.Cells.Item(132,11).Value = 11158.2729
$ws.Cells.Item(132,12).Value = 7045.5
$ws.Cells.Item(132,13).Value = -8628.2729
$ws.Cells.Item(132,14).Value = -12105.5
